$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pl_mw results for the 380 kV case (Case_5_238)
# Each inner array: row number, then new values for columns B, C, D, E, F, L, M, N
$newValues = @(
    @(2, 2.462849570000685, 0.2267414144939437, 0.1261810435928297, 0.05179948103840104, 2.512590215698481, 0.1935801646044268, 0.4435440823366008, 2.342519244637273),
    @(3, 2.33973469034612, 0.1971758706339131, 0.1267737299601492, 0.0512734343788015, 2.463575717351986, 0.1912509182811277, 0.4261271886266798, 2.353551735120377),
    @(4, 2.265725309365791, 0.1790397719284158, 0.1271821046594894, 0.05094355868138578, 2.43497697852662, 0.1899219057064343, 0.4157147108890342, 2.361022628853334),
    @(5, 2.235961576801287, 0.1716527671448773, 0.1273596603936547, 0.05080738384054229, 2.423696730668325, 0.1894057176093469, 0.4115419923478001, 2.364241640070574),
    @(6, 2.231043157541933, 0.1704263568157103, 0.127389815017164, 0.05078466603800624, 2.421846174276567, 0.1893215376320754, 0.4108533605601821, 2.364786677732639),
    @(7, 2.265322306450628, 0.178940134573395, 0.1271844541888676, 0.05094172927916851, 2.434823338025708, 0.1899148414408742, 0.4156581513566024, 2.361065335703088),
    @(8, 2.420069291190828, 0.2165431763658319, 0.1263761488004818, 0.0516195175781009, 2.495378032804183, 0.1927560197513642, 0.4374800784492265, 2.3461781948089),
    @(9, 2.736227738671687, 0.2904566404552327, 0.1251457441737998, 0.05289489744471965, 2.626112249226566, 0.1991326153030215, 0.4825258030473805, 2.322545106758824),
    @(10, 2.976467311397016, 0.3449235667933976, 0.1244606681742582, 0.05380036671954258, 2.729641756673601, 0.2043125403092318, 0.5170249657290782, 2.308615053884353),
    @(11, 3.08753569845743, 0.3697499892008977, 0.124197138447073, 0.05420570636084854, 2.778402473010146, 0.2067775197209869, 0.5330314114518799, 2.303032930292858),
    @(12, 3.129854200066745, 0.3791591112515675, 0.124104315524832, 0.05435827494959256, 2.797109177730334, 0.2077266295709137, 0.5391380859501851, 2.301028456949837),
    @(13, 3.120728571409188, 0.3771323240242737, 0.1241239957794136, 0.05432545739714101, 2.793069543975264, 0.2075215236522325, 0.5378208790629202, 2.301455279218928),
    @(14, 3.091012052373742, 0.3705239203229667, 0.1241893618623493, 0.05421827669387902, 2.779936617552011, 0.2068552890541469, 0.5335328992426227, 2.30286582463637),
    @(15, 3.072843688059493, 0.3664771351649847, 0.1242303097207866, 0.05415250555854279, 2.771923928000547, 0.2064492443251567, 0.5309123106331981, 2.303744090515366),
    @(16, 2.969244851402379, 0.3433021610678679, 0.1244788630098768, 0.0537737463771446, 2.726488856198415, 0.2041536384620599, 0.5159852336579789, 2.3089951162369),
    @(17, 2.906148922516934, 0.3290982369070434, 0.1246437018731257, 0.05353972241722715, 2.699044174495327, 0.2027732158910993, 0.5069083224728317, 2.312410338694363),
    @(18, 2.870025376367607, 0.3209330907955632, 0.1247430384959216, 0.05340449914525891, 2.683415237091339, 0.2019894521868935, 0.5017169206700842, 2.314445672688137),
    @(19, 2.857823252140463, 0.3181692734490014, 0.124777447912173, 0.05335860800106573, 2.678150347662353, 0.201725836619957, 0.4999642373630095, 2.315146969109108),
    @(20, 2.912848227931534, 0.3306097893930087, 0.1246256857373567, 0.05356469858216073, 2.70194948926121, 0.2029191062029412, 0.5078715295689662, 2.312039429459517),
    @(21, 3.099733452100338, 0.3724647477567373, 0.1241699726809955, 0.05424978320893281, 2.783787484600111, 0.207050552397277, 0.5347911477622631, 2.302448538720498),
    @(22, 3.223386689591337, 0.3998657928136709, 0.1239127892285552, 0.05469214454197591, 2.838685660932924, 0.2098420896092534, 0.5526493957191789, 2.296818138332512),
    @(23, 3.1572511737902, 0.3852368157880051, 0.1240463160833656, 0.05445653399992345, 2.809255337126046, 0.2083438108711704, 0.5430937461572043, 2.299764553755978),
    @(24, 2.909819001402241, 0.329926413837768, 0.1246338166059289, 0.05355340897678218, 2.700635531387746, 0.2028531184865585, 0.5074359794639633, 2.312206893705778),
    @(25, 2.649316820617571, 0.2704366529125934, 0.1254403567332787, 0.05255550338570814, 2.589444330953796, 0.197320958174025, 0.4700954177488015, 2.32833875154634)
)

$colIndexes = @(2, 3, 4, 5, 6, 12, 13, 14)

foreach ($rowArr in $newValues) {
    $r = [int]$rowArr[0]
    for ($i = 0; $i -lt $colIndexes.Count; $i++) {
        $ws.Cells.Item($r, [int]$colIndexes[$i]).Value = $rowArr[$i + 1]
    }
}
